# Refresh the "Updated symbol list" crypto snapshot: 2023-01-31 21:37 UTC run.
# Only cell text changes (prices, 1h volume %, and a few re-ranked coin rows/links) -
# no rows/columns are inserted or removed, so we just overwrite cell text in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price/Volume columns (D, E) must stay plain text (as in the source
# file, which stores them as inline strings) rather than become real numbers/percentages,
# so those are entered with a leading apostrophe - the normal Excel "store as text" prefix.

# Row 2
$ws.Range("D2").Value = "'311.87"
$ws.Range("E2").Value = "'1.82%"

# Row 3
$ws.Range("D3").Value = "'37.39"
$ws.Range("E3").Value = "'-0.18%"

# Row 4
$ws.Range("D4").Value = "'5.154"
$ws.Range("E4").Value = "'1.32%"

# Row 5
$ws.Range("D5").Value = "'0.07832"
$ws.Range("E5").Value = "'1.35%"

# Row 6
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.422"
$ws.Range("E6").Value = "'1.63%"

# Row 7
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'0.91%"

# Row 8
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.272"
$ws.Range("E8").Value = "'0.82%"

# Row 9
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.795"
$ws.Range("E9").Value = "'-7.38%"

# Row 10
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9183"
$ws.Range("E10").Value = "'0.00%"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1184"
$ws.Range("E11").Value = "'2.98%"

# Row 12
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1916"
$ws.Range("E12").Value = "'1.91%"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09093"
$ws.Range("E13").Value = "'4.13%"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03356"
$ws.Range("E14").Value = "'-1.41%"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09640"
$ws.Range("E15").Value = "'-0.69%"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001383"
$ws.Range("E16").Value = "'1.30%"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005762"
$ws.Range("E17").Value = "'-3.55%"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.510"
$ws.Range("E18").Value = "'-2.10%"

# Row 19
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'1.02%"

# Row 20
$ws.Range("D20").Value = "'5.259"
$ws.Range("E20").Value = "'4.74%"

# Row 21
$ws.Range("D21").Value = "'0.1273"
$ws.Range("E21").Value = "'-0.18%"

# Row 23
$ws.Range("D23").Value = "'0.04380"
$ws.Range("E23").Value = "'1.29%"

# Row 24
$ws.Range("D24").Value = "'0.001252"
$ws.Range("E24").Value = "'3.14%"

# Row 25
$ws.Range("D25").Value = "'0.004675"
$ws.Range("E25").Value = "'2.98%"

# Row 26
$ws.Range("D26").Value = "'0.0001363"
$ws.Range("E26").Value = "'0.78%"

# Row 27
$ws.Range("D27").Value = "'0.0004005"
$ws.Range("E27").Value = "'-98.10%"

# Row 39
$ws.Range("D39").Value = "'0.02285"
$ws.Range("E39").Value = "'3.54%"

# Row 40
$ws.Range("D40").Value = "'0.05051"
$ws.Range("E40").Value = "'2.69%"

# Row 41
$ws.Range("D41").Value = "'0.007435"
$ws.Range("E41").Value = "'-1.67%"

# Row 42
$ws.Range("D42").Value = "'0.009082"
$ws.Range("E42").Value = "'-8.70%"

# Row 43
$ws.Range("D43").Value = "'0.1354"
$ws.Range("E43").Value = "'1.45%"

# Row 44
$ws.Range("D44").Value = "'0.001956"
$ws.Range("E44").Value = "'-2.06%"

# Row 45
$ws.Range("D45").Value = "'0.009432"
$ws.Range("E45").Value = "'7.39%"

# Row 46
$ws.Range("D46").Value = "'0.00006688"
$ws.Range("E46").Value = "'2.26%"

# Row 47
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.20%"

# Row 48
$ws.Range("D48").Value = "'0.003297"
$ws.Range("E48").Value = "'9.91%"

# Row 49
$ws.Range("D49").Value = "'0.001004"
$ws.Range("E49").Value = "'-22.90%"

# Row 50
$ws.Range("D50").Value = "'0.00002108"
$ws.Range("E50").Value = "'0.20%"

# Row 51
$ws.Range("D51").Value = "'0.0002007"
$ws.Range("E51").Value = "'0.20%"
